$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (row 2 through 28) from 45212 to 45221 (date serial numbers)
$ws.Range("C2:C28").Value = 45221
